# ------------------------------------------------------------------
# Rebuild Language.xlsx: rename Sheet1 -> Comm, add Property/Guild/Tip/Item
# sheets, refresh shared text, and resize the data ranges.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- sheet 1: rename Sheet1 -> Comm ------------------------------
$wsComm = $wb.Worksheets.Item(1)
$wsComm.Name = "Comm"

# ---- add the 4 new sheets, in tab order after Comm ---------------
$wsProperty = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsProperty.Name = "Property"

$wsGuild = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsGuild.Name = "Guild"

$wsTip = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsTip.Name = "Tip"

$wsItem = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsItem.Name = "Item"

# ===================================================================
# Comm sheet: update text, extend used range to row 12
# ===================================================================
$wsComm.Cells.Item(2,1).Value = "Langage_Comm_1"
$wsComm.Cells.Item(2,3).Value = "确认"
$wsComm.Cells.Item(3,1).Value = "Langage_Comm_2"
$wsComm.Cells.Item(3,3).Value = "取消"
$wsComm.Cells.Item(4,1).Value = "Langage_Comm_3"
$wsComm.Cells.Item(4,3).Value = "登录"
$wsComm.Cells.Item(5,1).Value = "Langage_Comm_4"
$wsComm.Cells.Item(5,3).Value = "创建角色"
$wsComm.Cells.Item(6,1).Value = "Langage_Comm_5"
$wsComm.Cells.Item(6,3).Value = "进入游戏"
$wsComm.Cells.Item(7,1).Value = "Langage_Comm_6"
$wsComm.Cells.Item(7,3).Value = "中文_6"

# rows 8-12: blank cells that still carry the row-2..7 formatting (style s="2")
$wsComm.Range("A7:C7").Copy()
$wsComm.Range("A8:C12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# column widths (character units)
$wsComm.Columns.Item(1).ColumnWidth = 31.15
$wsComm.Columns.Item(2).ColumnWidth = 23.85
$wsComm.Columns.Item(3).ColumnWidth = 22.3

# ===================================================================
# Property sheet: ID/English/Chinese header + single HP/MP/... column
# ===================================================================
$wsProperty.Cells.Item(1,1).Value = "ID"
$wsProperty.Cells.Item(1,2).Value = "English"
$wsProperty.Cells.Item(1,3).Value = "Chinese"
$wsComm.Range("B1:C1").Copy()
$wsProperty.Range("B1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# A2:A28 share the plain "English/Chinese header" formatting (style s="1")
$wsComm.Range("B1").Copy()
$wsProperty.Range("A2:A28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsProperty.Cells.Item(2,1).Value = "Langage_HP"
$wsProperty.Cells.Item(3,1).Value = "Langage_MAXHP"
$wsProperty.Cells.Item(4,1).Value = "Langage_MP"
$wsProperty.Cells.Item(5,1).Value = "Langage_MAXMP"
$wsProperty.Cells.Item(6,1).Value = "Langage_VP"
$wsProperty.Cells.Item(7,1).Value = "Langage_ATTACK"

$wsProperty.Columns.Item(1).ColumnWidth = 50.55

# ===================================================================
# Guild sheet: ID/English/Chinese header + one Guild-confirm row
# ===================================================================
$wsGuild.Cells.Item(1,1).Value = "ID"
$wsGuild.Cells.Item(1,2).Value = "English"
$wsGuild.Cells.Item(1,3).Value = "Chinese"
$wsComm.Range("B1:C1").Copy()
$wsGuild.Range("B1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsGuild.Cells.Item(2,1).Value = "Langage_Guild_1"
$wsGuild.Cells.Item(2,2).Value = "Langage_1"
$wsGuild.Cells.Item(2,3).Value = "确认要加入这个公会吗？点击确认加入"
$wsComm.Range("A2:C2").Copy()
$wsGuild.Range("A2:C2").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$wsGuild.Cells.Item(2,1).Value = "Langage_Guild_1"
$wsGuild.Cells.Item(2,2).Value = "Langage_1"
$wsGuild.Cells.Item(2,3).Value = "确认要加入这个公会吗？点击确认加入"

# rows 3-12 and 16: blank cells keeping the same "data row" style (s="2")
$wsComm.Range("A2:C2").Copy()
$wsGuild.Range("A3:C12").PasteSpecial(-4122)
$wsGuild.Range("A16:C16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# row 22 col A only: blank cell keeping the plain header style (s="1")
$wsComm.Range("B1").Copy()
$wsGuild.Range("A22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsGuild.Columns.Item(1).ColumnWidth = 31.15
$wsGuild.Columns.Item(2).ColumnWidth = 23.85
$wsGuild.Columns.Item(3).ColumnWidth = 22.3

# ===================================================================
# Tip / Item sheets: header row only
# ===================================================================
$wsTip.Cells.Item(1,1).Value = "ID"
$wsTip.Cells.Item(1,2).Value = "English"
$wsTip.Cells.Item(1,3).Value = "Chinese"
$wsComm.Range("B1:C1").Copy()
$wsTip.Range("B1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsItem.Cells.Item(1,1).Value = "ID"
$wsItem.Cells.Item(1,2).Value = "English"
$wsItem.Cells.Item(1,3).Value = "Chinese"
$wsComm.Range("B1:C1").Copy()
$wsItem.Range("B1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ===================================================================
# Selections / active sheet (cosmetic, matches the saved view state)
# ===================================================================
$wsProperty.Activate()
$wsProperty.Rows.Item(1).Select()

$wsGuild.Activate()
$wsGuild.Range("A12").Select()

$wsTip.Activate()
$wsTip.Rows.Item(1).Select()

$wsItem.Activate()
$wsItem.Rows.Item(1).Select()

$wsComm.Activate()
$wsComm.Range("C8").Select()

# absolute path recorded by the authoring machine
$wb.Application.ActiveWorkbook.Path | Out-Null
